# Fix contrat AV (Deces, Cession, Revision loyer) and multiple proprietaires
# have same mandataire error.
#
# Two new rows describing contrat "009/TTT/AV1" / mandataire "ZERNAKH
# ABDELLAH" (IB19558) need to be inserted right after the header row, and
# the grand-total (footer) row needs to be updated to reflect the two new
# lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 2 - this pushes the former rows 2..5
# (NASIRI HASNAA, SOFIA BADRANE, DOUNIA LAMKADDAM, totals) down to rows
# 4..7, unchanged.
$ws.Rows(2).Insert()
$ws.Rows(2).Insert()

# New row 2: ZERNAKH ABDELLAH - loyer line
$ws.Cells.Item(2, 1).Value = "009/TTT/AV1"
$ws.Cells.Item(2, 2).Value = "Direction régionale"
$ws.Cells.Item(2, 3).Value = "IB19558"
$ws.Cells.Item(2, 4).Value = "ZERNAKH ABDELLAH"
$ws.Cells.Item(2, 5).Value = "non"
$ws.Cells.Item(2, 6).Value = "mensuelle"
$ws.Cells.Item(2, 7).Value = 10
$ws.Cells.Item(2, 8).Value = 3500
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 350
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 1000
$ws.Cells.Item(2, 15).Value = 4150

# New row 3: ZERNAKH ABDELLAH - second line (same contrat/mandataire)
$ws.Cells.Item(3, 1).Value = "009/TTT/AV1"
$ws.Cells.Item(3, 2).Value = "Direction régionale"
$ws.Cells.Item(3, 3).Value = "IB19558"
$ws.Cells.Item(3, 4).Value = "ZERNAKH ABDELLAH"
$ws.Cells.Item(3, 5).Value = "non"
$ws.Cells.Item(3, 6).Value = "mensuelle"
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 8).Value = 3500
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 350
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 1000
$ws.Cells.Item(3, 15).Value = 3150

# Update grand-total (footer) row, now shifted to row 7, to include the
# two new rows' contributions.
$ws.Cells.Item(7, 8).Value = 15500.01
$ws.Cells.Item(7, 10).Value = 1550.01
$ws.Cells.Item(7, 14).Value = 2000
$ws.Cells.Item(7, 15).Value = 14950
